# "Con cambios y errores" - update the product names (column C) and turn
# the Description (column B) into a computed "desc" + name formula; also
# add a new barcode-ish concat column (I) and fix page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product names for rows 2-5 (name / default_code column C)
$newNames = @{
    2 = "Oscar8"
    3 = "Irving8"
    4 = "Chido8"
    5 = "OCHO"
}

foreach ($row in 2..5) {
    $name = $newNames[$row]

    # C: the product name itself
    $ws.Cells.Item($row, 3).Value = $name

    # B: Description becomes a formula built from the (new) name in C
    $ws.Cells.Item($row, 2).Formula = "=CONCAT(`"desc`",C$row)"

    # I: new column concatenating A, B and C together
    $ws.Cells.Item($row, 9).Formula = "=CONCAT(A$row,B$row,C$row)"
}

# Switch the page to portrait orientation
$ws.PageSetup.Orientation = 1

# Leave the selection on B5, like in the edited workbook
$ws.Range("B5").Select()
